$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "food super group"
$ws.Range("B7").Value = "dairy"

$ws.Range("A8").Select()
